$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 0.0002034219952677446
$ws.Range("E2").Value = 0.0002034219952677446

# Row 3
$ws.Range("D3").Value = 0.9999998010092973
$ws.Range("E3").Value = 0.9999998010092973

# Row 4
$ws.Range("D4").Value = 0.3283828948148235
$ws.Range("E4").Value = 0.3283828948148235

# Row 5
$ws.Range("D5").Value = [double]"6.09921176195907E-27"
$ws.Range("E5").Value = [double]"6.09921176195907E-27"

# Row 6
$ws.Range("D6").Value = 0.1562520927195231
$ws.Range("E6").Value = 0.1562520927195231

# Row 7
$ws.Range("D7").Value = 0.9999999996625948
$ws.Range("E7").Value = [double]"3.374052148785722E-10"

# Row 8
$ws.Range("D8").Value = 0.9999999966402207
$ws.Range("E8").Value = [double]"3.359779343625746E-09"

# Row 9
$ws.Range("D9").Value = 0.0001213909939871626
$ws.Range("E9").Value = 0.9998786090060129

# Row 11
$ws.Range("D11").Value = [double]"3.561512789112387E-07"
$ws.Range("E11").Value = 0.9999996438487211
$ws.Range("F11").Value = 3.986258268356323
